# Updated cryptos list (Price / Volume(1h) columns) per the commit diff.
# Each entry: spreadsheet Row, new Price (column D, `$null` when unchanged),
# and new Volume(1h) (column E, always present in this diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = '26.330.70'; Volume = '  +1.11%  ' },
    @{ Row = 3; Price = '1.682.26'; Volume = '  +0.86%  ' },
    @{ Row = 4; Price = $null; Volume = '  +0.42%  ' },
    @{ Row = 5; Price = '218.17'; Volume = '  +0.59%  ' },
    @{ Row = 6; Price = '0.5505'; Volume = '  +7.90%  ' },
    @{ Row = 7; Price = $null; Volume = '  +0.35%  ' },
    @{ Row = 8; Price = $null; Volume = '  +1.72%  ' },
    @{ Row = 9; Price = '0.06488'; Volume = '  +1.24%  ' },
    @{ Row = 10; Price = '22.07'; Volume = '  +0.96%  ' },
    @{ Row = 11; Price = '0.07549'; Volume = '  +1.50%  ' },
    @{ Row = 12; Price = '4.544'; Volume = '  +0.82%  ' },
    @{ Row = 13; Price = '1.675.73'; Volume = '  +0.38%  ' },
    @{ Row = 14; Price = '0.5810'; Volume = '  -0.57%  ' },
    @{ Row = 15; Price = '0.000008433'; Volume = '  -1.77%  ' },
    @{ Row = 16; Price = '64.98'; Volume = '  +0.92%  ' },
    @{ Row = 17; Price = '26.349.24'; Volume = '  +1.07%  ' },
    @{ Row = 18; Price = '4.931'; Volume = '  -0.31%  ' },
    @{ Row = 19; Price = $null; Volume = '  +0.35%  ' },
    @{ Row = 20; Price = '10.93'; Volume = '  +1.27%  ' },
    @{ Row = 21; Price = '191.10'; Volume = '  -0.47%  ' },
    @{ Row = 22; Price = '6.229'; Volume = '  +0.29%  ' },
    @{ Row = 23; Price = '1.008'; Volume = '  +0.36%  ' },
    @{ Row = 24; Price = '146.89'; Volume = '  +1.50%  ' },
    @{ Row = 25; Price = '0.1314'; Volume = '  +9.66%  ' },
    @{ Row = 26; Price = '7.908'; Volume = '  +3.71%  ' },
    @{ Row = 27; Price = '15.81'; Volume = '  +0.81%  ' },
    @{ Row = 28; Price = '0.06324'; Volume = '  -2.29%  ' },
    @{ Row = 29; Price = '1.391'; Volume = '  +4.84%  ' },
    @{ Row = 30; Price = '1.323'; Volume = '  +0.58%  ' },
    @{ Row = 31; Price = '3.592'; Volume = '  +1.33%  ' },
    @{ Row = 32; Price = '3.580'; Volume = '  +1.70%  ' },
    @{ Row = 33; Price = '1.666'; Volume = '  +0.99%  ' },
    @{ Row = 34; Price = $null; Volume = '  +1.58%  ' },
    @{ Row = 35; Price = '0.6191'; Volume = '  +1.37%  ' },
    @{ Row = 36; Price = $null; Volume = '  +1.45%  ' },
    @{ Row = 37; Price = '2.717'; Volume = '  +1.37%  ' },
    @{ Row = 38; Price = '6.233'; Volume = '  -0.62%  ' },
    @{ Row = 39; Price = '1.112.42'; Volume = '  +1.94%  ' },
    @{ Row = 40; Price = $null; Volume = '  +1.26%  ' },
    @{ Row = 41; Price = '0.8715'; Volume = '  +0.83%  ' },
    @{ Row = 42; Price = '1.016'; Volume = '  +0.74%  ' },
    @{ Row = 43; Price = '100.75'; Volume = '  -0.12%  ' },
    @{ Row = 44; Price = '1.831.97'; Volume = '  +0.84%  ' },
    @{ Row = 45; Price = $null; Volume = '  -5.16%  ' },
    @{ Row = 46; Price = $null; Volume = '  +1.53%  ' },
    @{ Row = 47; Price = '8.199'; Volume = '  +1.91%  ' },
    @{ Row = 48; Price = '1.005'; Volume = '  -0.11%  ' },
    @{ Row = 49; Price = '0.05275'; Volume = '  +0.84%  ' },
    @{ Row = 50; Price = '0.4293'; Volume = '  +0.26%  ' },
    @{ Row = 51; Price = '6.064'; Volume = '  +0.15%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Cells.Item($u.Row, 4)
        # Force Text storage for price strings that would otherwise be
        # auto-coerced to a number/date by Excel (e.g. '218.17', '0.5505').
        # Values that already contain 2+ dots (e.g. '26.330.70') are never
        # parseable as a number, so they stay Text without this, and we
        # skip it there to avoid an unnecessary style change.
        $dotCount = ([regex]::Matches($u.Price, '\.')).Count
        if ($dotCount -lt 2) {
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $u.Price
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.Volume
}
